$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a serial date value that was bumped by one day
# (46082 -> 46083) for every data row (rows 2 through 550).
$ws.Range("C2:C550").Value = 46083
